# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity, Temperature and mmWave
# sheets, matching the latest batch of readings captured on 2026-01-28.
#
# All the log columns (Date, Timestamp, Hour, Value, ...) are stored as
# plain text in this workbook, even though several of them look like
# dates/times/percentages. Setting a Range.Value directly would let Excel
# "helpfully" reinterpret e.g. "2026-01-28" as a date serial or "88.0%" as
# a numeric percentage (and stick a NumberFormat on the cell in the
# process). To avoid that, force each new row's format to Text ("@")
# before writing the values, then clear the formatting again afterwards
# so the cells end up back at the sheet's default (General) style with
# plain text content - matching every other row already in the log.
#
# NOTE: named parameters (-Foo bar) aren't reliable in this PowerShell
# engine, so Add-LogRow takes its arguments positionally.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($Sheet, $Row, $Date, $Timestamp, $Hour, $Location, $Value, $Status)

    $rng = $Sheet.Range($Sheet.Cells.Item($Row, 1), $Sheet.Cells.Item($Row, 6))
    $rng.NumberFormat = "@"

    $Sheet.Cells.Item($Row, 1).Value = $Date
    $Sheet.Cells.Item($Row, 2).Value = $Timestamp
    $Sheet.Cells.Item($Row, 3).Value = $Hour
    $Sheet.Cells.Item($Row, 4).Value = $Location
    $Sheet.Cells.Item($Row, 5).Value = $Value
    $Sheet.Cells.Item($Row, 6).Value = $Status

    $rng.ClearFormats()
}

# ---------------------------------------------------------------------
# PIR sheet: append rows 74-78
# ---------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")

Add-LogRow $pir 74 "2026-01-28" "16:41:14" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $pir 75 "2026-01-28" "16:41:14" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $pir 76 "2026-01-28" "16:41:18" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $pir 77 "2026-01-28" "16:41:23" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $pir 78 "2026-01-28" "16:41:28" "16:00" "Bathroom" "No Motion" "Inactive"

# ---------------------------------------------------------------------
# Humidity sheet: append rows 73-76
# ---------------------------------------------------------------------
$humidity = $wb.Worksheets.Item("Humidity")

Add-LogRow $humidity 73 "2026-01-28" "16:41:13" "16:00" "Bathroom" "87.1%" "Active"
Add-LogRow $humidity 74 "2026-01-28" "16:41:15" "16:00" "Bathroom" "88.0%" "Active"
Add-LogRow $humidity 75 "2026-01-28" "16:41:24" "16:00" "Bathroom" "88.0%" "Active"
Add-LogRow $humidity 76 "2026-01-28" "16:41:28" "16:00" "Bathroom" "87.2%" "Active"

# ---------------------------------------------------------------------
# Temperature sheet: append rows 73-76
# ---------------------------------------------------------------------
$temperature = $wb.Worksheets.Item("Temperature")

Add-LogRow $temperature 73 "2026-01-28" "16:41:14" "16:00" "Bathroom" "22.8C" "Active"
Add-LogRow $temperature 74 "2026-01-28" "16:41:16" "16:00" "Bathroom" "22.8C" "Active"
Add-LogRow $temperature 75 "2026-01-28" "16:41:24" "16:00" "Bathroom" "22.8C" "Active"
Add-LogRow $temperature 76 "2026-01-28" "16:41:28" "16:00" "Bathroom" "22.8C" "Active"

# ---------------------------------------------------------------------
# mmWave sheet: append row 2
# ---------------------------------------------------------------------
$mmwave = $wb.Worksheets.Item("mmWave")

Add-LogRow $mmwave 2 "2026-01-28" "16:41:14" "16:00" "Living Room" "NO_PRESENCE" "Active"
